# daily auto push: 2026-01-27 22:35 UTC
#
# A new sample row is inserted into the "sei1" daily log sheet at row 737
# (date 2026/01/28, weekday 水, time 3, ranking 201). Every row that used
# to live at 737..778 shifts down by one (738..779), and the sheet's
# dimension grows from A1:D778 to A1:D779.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 737..778 down one slot by inserting a blank row at 737 (Excel
# shifts everything below it down automatically, formatting included).
$ws.Rows("737:737").Insert()

# Column A holds dates formatted as plain text (e.g. "2026/01/27"), not
# real date serials. Force the new cell to Text first so the COM layer
# doesn't auto-parse the "YYYY/MM/DD" string into a date value, then
# restore the Normal style so no stray number-format sticks around.
$ws.Range("A737").NumberFormat = "@"
$ws.Range("A737").Value = "2026/01/28"
$ws.Range("A737").Style = "Normal"

$ws.Range("B737").Value = "水"
$ws.Range("C737").Value = 3
$ws.Range("D737").Value = 201
